# Remove the "Pan" (C) and "Primary Email" (D) investor-creation columns
# from the Portfolio Investments import sheet, along with the hyperlinks
# (mailto: links on the Primary Email column) that went with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto hyperlinks that live on the Primary Email column (D2/D5/D8)
# before the columns shift, so no stray/mis-targeted hyperlinks remain.
$ws.Hyperlinks.Delete()

# Select columns C:D first so the post-delete selection matches what Excel
# leaves behind (the columns that slide into the deleted ones stay selected).
$ws.Range("C:D").Select() | Out-Null

# Delete columns C (Pan) and D (Primary Email) entirely - this shifts every
# later column left by two (old E -> C, old F -> D, etc.).
$ws.Range("C:D").Delete()
